$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.817.24'
$ws.Range('E2').Value = '  +0.01%  '
$ws.Range('D3').Value = '1.542.40'
$ws.Range('E3').Value = '  -1.58%  '
$ws.Range('E4').Value = '  +0.26%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range('D5').Value = '206.01'
$ws.Range('E5').Value = '  -0.21%  '
$ws.Range('E6').Value = '  -0.70%  '
$ws.Range('E8').Value = '  -0.47%  '
$ws.Range('E9').Value = '  -2.64%  '
$ws.Range('E10').Value = '  -0.57%  '
$ws.Range('E11').Value = '  -1.10%  '
$ws.Range('D12').Value = '1.762.61'
$ws.Range('E12').Value = '  -1.48%  '
$ws.Range('D13').Value = '1.543.80'
$ws.Range('E13').Value = '  -1.25%  '
$ws.Range('E14').Value = '  -1.51%  '
$ws.Range('E15').Value = '  -1.05%  '
$ws.Range('D16').Value = '26.821.82'
$ws.Range('E16').Value = '  +0.03%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range('D17').Value = '61.26'
$ws.Range('E17').Value = '  -0.31%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range('D18').Value = '214.98'
$ws.Range('E18').Value = '  +0.12%  '
$ws.Range('E19').Value = '  -2.62%  '
$ws.Range('E21').Value = '  +0.21%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range('D22').Value = '4.00'
$ws.Range('E22').Value = '  -3.04%  '
$ws.Range('E23').Value = '  -1.36%  '
$ws.Range('E24').Value = '  -2.82%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range('D25').Value = '153.03'
$ws.Range('E25').Value = '  -0.32%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range('D26').Value = '6.61'
$ws.Range('E26').Value = '  -1.98%  '
$ws.Range('E27').Value = '  -0.98%  '
$ws.Range('E28').Value = '  +0.20%  '
$ws.Range('E29').Value = '  -0.65%  '
$ws.Range('E30').Value = '  -1.87%  '
$ws.Range('E31').Value = '  -1.50%  '
$ws.Range('E32').Value = '  +1.56%  '
$ws.Range('D33').Value = '1.369.33'
$ws.Range('E33').Value = '  -1.93%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range('D34').Value = '2.95'
$ws.Range('E34').Value = '  +0.67%  '
$ws.Range('E35').Value = '  -1.26%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range('D36').Value = '0.963'
$ws.Range('E36').Value = '  +2.87%  '
$ws.Range('E37').Value = '  +0.02%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range('D38').Value = '0.0165'
$ws.Range('E38').Value = '  +1.46%  '
$ws.Range('E39').Value = '  -1.45%  '
$ws.Range('E40').Value = '  +9.07%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range('D41').Value = '0.807'
$ws.Range('E41').Value = '  -1.11%  '
$ws.Range('B42').Value = 'WEMIXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D42").NumberFormat = "@"
$ws.Range('D42').Value = '0.989'
$ws.Range('E42').Value = '  -0.13%  '
$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D43").NumberFormat = "@"
$ws.Range('D43').Value = '2.21'
$ws.Range('E43').Value = '  +1.43%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").NumberFormat = "@"
$ws.Range('D44').Value = '63.20'
$ws.Range('E44').Value = '  -0.22%  '
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D45").NumberFormat = "@"
$ws.Range('D45').Value = '1.74'
$ws.Range('E45').Value = '  -3.46%  '
$ws.Range('B46').Value = 'RocketPoolETH'
$ws.Range('C46').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D46').Value = '1.677.05'
$ws.Range('E46').Value = '  -1.50%  '
$ws.Range('B47').Value = 'Quant'
$ws.Range('C47').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range('D47').Value = '84.17'
$ws.Range('E47').Value = '  -2.19%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").NumberFormat = "@"
$ws.Range('D48').Value = '0.0510'
$ws.Range('E48').Value = '  +3.64%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.0₇0970'
$ws.Range('E49').Value = '  -1.57%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").NumberFormat = "@"
$ws.Range('D50').Value = '0.0941'
$ws.Range('E50').Value = '  -1.32%  '
$ws.Range('E51').Value = '  -0.03%  '
